$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = "57.304.17"
$ws.Cells.Item(2,5).Value = "  -0.35%  "

# Row 3
$ws.Cells.Item(3,4).Value = "2.358.71"
$ws.Cells.Item(3,5).Value = "  +0.99%  "

# Row 4
$ws.Cells.Item(4,5).Value = "  -0.12%  "

# Row 5
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "519.83"
$ws.Cells.Item(5,5).Value = "  -0.26%  "

# Row 6
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "135.56"
$ws.Cells.Item(6,5).Value = "  +0.33%  "

# Row 7
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = "0.996"
$ws.Cells.Item(7,5).Value = "  -0.10%  "

# Row 8
$ws.Cells.Item(8,5).Value = "  -0.03%  "

# Row 9
$ws.Cells.Item(9,5).Value = "  -1.35%  "

# Row 10
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = "5.48"
$ws.Cells.Item(10,5).Value = "  +4.88%  "

# Row 11
$ws.Cells.Item(11,5).Value = "  -0.69%  "

# Row 12
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = "0.343"
$ws.Cells.Item(12,5).Value = "  -0.67%  "

# Row 13
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = "24.37"
$ws.Cells.Item(13,5).Value = "  +1.56%  "

# Row 14
$ws.Cells.Item(14,4).Value = "2.778.00"
$ws.Cells.Item(14,5).Value = "  +0.71%  "

# Row 15
$ws.Cells.Item(15,4).Value = "57.292.34"
$ws.Cells.Item(15,5).Value = "  +0.13%  "

# Row 16
$ws.Cells.Item(16,5).Value = "  -0.64%  "

# Row 17
$ws.Cells.Item(17,4).Value = "2.368.40"
$ws.Cells.Item(17,5).Value = "  +0.33%  "

# Row 18
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "10.58"
$ws.Cells.Item(18,5).Value = "  -0.31%  "

# Row 19
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "328.92"
$ws.Cells.Item(19,5).Value = "  +1.87%  "

# Row 20
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = "4.25"
$ws.Cells.Item(20,5).Value = "  -1.37%  "

# Row 21
$ws.Cells.Item(21,5).Value = "  +0.01%  "

# Row 22
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "0.999"
$ws.Cells.Item(22,5).Value = "  -0.01%  "

# Row 23
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "61.38"
$ws.Cells.Item(23,5).Value = "  -0.46%  "

# Row 24
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "8.92"
$ws.Cells.Item(24,5).Value = "  +13.93%  "

# Row 25
$ws.Cells.Item(25,5).Value = "  +2.13%  "

# Row 26
$ws.Cells.Item(26,5).Value = "  +0.70%  "

# Row 27
$ws.Cells.Item(27,5).Value = "  +10.17%  "

# Row 28
$ws.Cells.Item(28,4).Value = "0.0₃0744"
$ws.Cells.Item(28,5).Value = "  -0.12%  "

# Row 29
$ws.Cells.Item(29,2).Value = "PancakeSwap"
$ws.Cells.Item(29,3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "1.70"
$ws.Cells.Item(29,5).Value = "  +0.20%  "

# Row 30
$ws.Cells.Item(30,2).Value = "Monero"
$ws.Cells.Item(30,3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "167.03"
$ws.Cells.Item(30,5).Value = "  -2.87%  "

# Row 31
$ws.Cells.Item(31,5).Value = "  -0.75%  "

# Row 32
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = "18.59"
$ws.Cells.Item(32,5).Value = "  +0.61%  "

# Row 33
$ws.Cells.Item(33,5).Value = "  +0.06%  "

# Row 34
$ws.Cells.Item(34,5).Value = "  +2.05%  "

# Row 35
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = "0.994"
$ws.Cells.Item(35,5).Value = "  -0.28%  "

# Row 36
$ws.Cells.Item(36,5).Value = "  -3.59%  "

# Row 37
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = "4.03"
$ws.Cells.Item(37,5).Value = "  -0.84%  "

# Row 38
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "1.61"
$ws.Cells.Item(38,5).Value = "  +4.95%  "

# Row 39
$ws.Cells.Item(39,5).Value = "  +3.13%  "

# Row 40
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "149.71"
$ws.Cells.Item(40,5).Value = "  +6.61%  "

# Row 41
$ws.Cells.Item(41,5).Value = "  +0.53%  "

# Row 42
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "3.65"
$ws.Cells.Item(42,5).Value = "  +1.01%  "

# Row 43
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "288.37"
$ws.Cells.Item(43,5).Value = "  +2.35%  "

# Row 44
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "5.29"
$ws.Cells.Item(44,5).Value = "  +1.68%  "

# Row 45
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "0.0939"
$ws.Cells.Item(45,5).Value = "  +0.77%  "

# Row 46
$ws.Cells.Item(46,5).Value = "  -0.65%  "

# Row 47
$ws.Cells.Item(47,5).Value = "  +0.02%  "

# Row 48
$ws.Cells.Item(48,5).Value = "  +4.80%  "

# Row 49
$ws.Cells.Item(49,5).Value = "  +1.14%  "

# Row 50
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "17.70"
$ws.Cells.Item(50,5).Value = "  +3.62%  "

# Row 51
$ws.Cells.Item(51,2).Value = "WhiteBITCoin"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "10.98"
$ws.Cells.Item(51,5).Value = "  +1.08%  "
